# Adds the "SAMPLE TEST REPORT (from canvas)" outline content as new
# paragraphs at the end of the document body, just before the final
# (pre-existing) empty paragraph / sectPr.
#
# We build each new paragraph as a small OOXML fragment and insert it
# via Range.InsertXML(WordOpenXML-package). InsertXML is called once per
# paragraph (rather than once for the whole block) because this host's
# InsertXML implementation replaces the *entire* body instead of doing a
# true range-insert when the fragment contains more than one top-level
# <w:p>.

$d = $word.ActiveDocument

function Insert-WordXmlParagraph {
    param(
        [System.__ComObject]$TargetRange,
        [string]$ParagraphXml
    )
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $ParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$TargetRange.InsertXML($pkg)
}

$newParagraphs = @(
@'
<w:p><w:r><w:t>SAMPLE TEST REPORT (from canvas)</w:t></w:r></w:p>
'@,
@'
<w:p><w:r><w:t>CS3300 PROJECT NAME TEST REPORT</w:t></w:r><w:r><w:br/><w:t>DATE</w:t></w:r><w:r><w:br/><w:t>Report covers all testing, test schedule, product features and quality, a list of all flags remaining</w:t></w:r><w:r><w:br/><w:t>open with a RISK assessment, and lessons learned.</w:t></w:r><w:r><w:br/><w:t>Test Report Outline</w:t></w:r><w:r><w:br/><w:t>1. Project Title:</w:t></w:r><w:r><w:br/><w:t>2. Test Report Scope: Report on the increments 1, 2, 3 and 4 testing associated with release</w:t></w:r><w:r><w:br/><w:t>1 and 2 testing, and system testing.</w:t></w:r><w:r><w:br/><w:t>3. Test team: QA/Test lead goes first</w:t></w:r><w:r><w:br/><w:t>4. Test Schedule:</w:t></w:r><w:r><w:br/><w:t>a. Unit/Increment testing occurred 17 November, with revision testing on 20</w:t></w:r><w:r><w:br/><w:t>November</w:t></w:r><w:r><w:br/><w:t>b. Component/Release testing occurred 21 November, with revision testing on 23</w:t></w:r><w:r><w:br/><w:t>November</w:t></w:r><w:r><w:br/><w:t>c. System testing occurred 1 December</w:t></w:r><w:r><w:br/><w:t>5. Test Description: describe each test, what was tested and the test outcome.</w:t></w:r><w:r><w:br/><w:t>6. Final product features and quality: describe the final product and any flaws</w:t></w:r><w:r><w:br/><w:t>7. Test Flags and Risk Management: List all open flaws/flags remaining and provide a</w:t></w:r><w:r><w:br/><w:t>RISK assessment for each</w:t></w:r><w:r><w:br/><w:t>8. Lessons Learned from testing</w:t></w:r></w:p>
'@,
@'
<w:p/>
'@
)

foreach ($paraXml in $newParagraphs) {
    # Always re-fetch the last paragraph: each insert lands right before
    # it, so the (originally trailing, still-empty) final paragraph keeps
    # getting pushed later - exactly mirroring "insert before the closing
    # empty paragraph" semantics.
    $anchor = $d.Paragraphs($d.Paragraphs.Count).Range
    Insert-WordXmlParagraph $anchor $paraXml
}
